$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 302.46155
$ws.Range("I9").Value = 266.54544
$ws.Range("K9").Value = 266.54544
$ws.Range("M9").Value = -97.54543999999999
$ws.Range("H43").Value = 1214.0476
$ws.Range("I43").Value = 920.25
$ws.Range("J43").Value = 1283.1765
$ws.Range("K43").Value = 920.25
$ws.Range("L43").Value = 1283.1765
$ws.Range("M43").Value = -851.25
$ws.Range("N43").Value = -1421.1765
$ws.Range("H51").Value = 6112.5
$ws.Range("H62").Value = 3017.182
$ws.Range("I62").Value = 1658
$ws.Range("K62").Value = 1658
$ws.Range("M62").Value = -1034
$ws.Range("H65").Value = 3017.182
$ws.Range("I65").Value = 1658
$ws.Range("K65").Value = 8290
$ws.Range("M65").Value = -5170
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H113").Value = 11842.5
$ws.Range("I113").Value = 6501.3335
$ws.Range("J113").Value = 13622.889
$ws.Range("K113").Value = 6501.3335
$ws.Range("L113").Value = 13622.889
$ws.Range("M113").Value = -3247.3335
$ws.Range("N113").Value = -20130.889
$ws.Range("H132").Value = 131593.06
$ws.Range("I132").Value = 144620.67
$ws.Range("K132").Value = 433862.01
$ws.Range("M132").Value = -431332.01

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2100
$ws.Range("I2").Value = 2316.6667
$ws.Range("J2").Value = 1666.6666
$ws.Range("K2").Value = 2316.6667
$ws.Range("L2").Value = 1666.6666
$ws.Range("M2").Value = -2203.6667
$ws.Range("N2").Value = -1892.6666
$ws.Range("H32").Value = 6305.148
$ws.Range("I32").Value = 4131.6177
$ws.Range("J32").Value = 10000.15
$ws.Range("K32").Value = 4131.6177
$ws.Range("L32").Value = 10000.15
$ws.Range("M32").Value = -3844.6177
$ws.Range("N32").Value = -10574.15
$ws.Range("H61").Value = 1661.037
$ws.Range("I61").Value = 1352
$ws.Range("K61").Value = 1352
$ws.Range("M61").Value = -1140
$ws.Range("H116").Value = 2100
$ws.Range("I116").Value = 2316.6667
$ws.Range("J116").Value = 1666.6666
$ws.Range("K116").Value = 2316.6667
$ws.Range("L116").Value = 1666.6666
$ws.Range("M116").Value = -22.66670000000022
$ws.Range("N116").Value = -6254.6666
$ws.Range("H132").Value = 2281.9805
$ws.Range("I132").Value = 1586.4
$ws.Range("J132").Value = 4811.364
$ws.Range("K132").Value = 4759.200000000001
$ws.Range("L132").Value = 14434.092
$ws.Range("M132").Value = -2229.200000000001
$ws.Range("N132").Value = -19494.092
$ws.Range("H136").Value = 1661.037
$ws.Range("I136").Value = 1352
$ws.Range("K136").Value = 4056
$ws.Range("M136").Value = -1506

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2100
$ws.Range("I3").Value = 2316.6667
$ws.Range("J3").Value = 1666.6666
$ws.Range("K3").Value = 2316.6667
$ws.Range("L3").Value = 1666.6666
$ws.Range("M3").Value = -2202.6667
$ws.Range("N3").Value = -1894.6666
$ws.Range("H94").Value = 966.2222
$ws.Range("I94").Value = 819.6786
$ws.Range("K94").Value = 819.6786
$ws.Range("M94").Value = -368.6786
$ws.Range("H99").Value = 1709.5938
$ws.Range("I99").Value = 1232.762
$ws.Range("J99").Value = 2619.9092
$ws.Range("K99").Value = 1232.762
$ws.Range("L99").Value = 2619.9092
$ws.Range("M99").Value = 265.2380000000001
$ws.Range("N99").Value = -5615.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 28749.75
$ws.Range("I21").Value = 35000
$ws.Range("J21").Value = 26666.334
$ws.Range("K21").Value = 35000
$ws.Range("L21").Value = 26666.334
$ws.Range("M21").Value = -34765
$ws.Range("N21").Value = -27136.334
$ws.Range("H32").Value = 7140
$ws.Range("I32").Value = 710
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 710
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -394
$ws.Range("N32").Value = -20632
$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -21186
$ws.Range("H132").Value = 3419.1714
$ws.Range("I132").Value = 3529.611
$ws.Range("K132").Value = 10588.833
$ws.Range("M132").Value = -8058.832999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1299
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 1299
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 3897
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -5015
$ws.Range("H131").Value = 7576733.5
$ws.Range("I131").Value = 83335220
$ws.Range("J131").Value = 885.65
$ws.Range("K131").Value = 250005660
$ws.Range("L131").Value = 2656.95
$ws.Range("M131").Value = -250000620
$ws.Range("N131").Value = -12736.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 9260035
$ws.Range("I107").Value = 767
$ws.Range("J107").Value = 18519302
$ws.Range("K107").Value = 767
$ws.Range("L107").Value = 18519302
$ws.Range("M107").Value = 1153
$ws.Range("N107").Value = -18523142
$ws.Range("H116").Value = 29000
$ws.Range("J116").Value = 29000
$ws.Range("L116").Value = 29000
$ws.Range("N116").Value = -38178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3978.6428
$ws.Range("I7").Value = 2400.389
$ws.Range("J7").Value = 6819.5
$ws.Range("K7").Value = 2400.389
$ws.Range("L7").Value = 6819.5
$ws.Range("M7").Value = -2288.389
$ws.Range("N7").Value = -7043.5
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -683
$ws.Range("H55").Value = 395.29413
$ws.Range("I55").Value = 243.25
$ws.Range("J55").Value = 530.44446
$ws.Range("K55").Value = 243.25
$ws.Range("L55").Value = 530.44446
$ws.Range("M55").Value = -70.25
$ws.Range("N55").Value = -876.44446
$ws.Range("H126").Value = 3978.6428
$ws.Range("I126").Value = 2400.389
$ws.Range("J126").Value = 6819.5
$ws.Range("K126").Value = 7201.167
$ws.Range("L126").Value = 20458.5
$ws.Range("M126").Value = -4731.167
$ws.Range("N126").Value = -25398.5
$ws.Range("H136").Value = 2639.2632
$ws.Range("I136").Value = 1451.4193
$ws.Range("J136").Value = 7899.7144
$ws.Range("K136").Value = 4354.257900000001
$ws.Range("L136").Value = 23699.1432
$ws.Range("M136").Value = -1804.257900000001
$ws.Range("N136").Value = -28799.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2063.262
$ws.Range("I126").Value = 1429.6364
$ws.Range("J126").Value = 2760.25
$ws.Range("K126").Value = 4288.9092
$ws.Range("L126").Value = 8280.75
$ws.Range("M126").Value = -1818.9092
$ws.Range("N126").Value = -13220.75
$ws.Range("H136").Value = 1002.125
$ws.Range("I136").Value = 547.129
$ws.Range("J136").Value = 1831.8235
$ws.Range("K136").Value = 1641.387
$ws.Range("L136").Value = 5495.470499999999
$ws.Range("M136").Value = 908.6129999999998
$ws.Range("N136").Value = -10595.4705

Write-Output "Edit complete"